# Scale the table's columns proportionally (duplicate-column style
# resize): keep the overall table width the same, but redistribute it
# across the 3 columns with new proportions.
#
# EMU widths (before -> after), total stays 8444091 EMU:
#   col 1: 2814697 -> 1629100
#   col 2: 2814697 -> 5112912
#   col 3: 2814697 -> 1702079
#
# PowerPoint's COM object model reports/accepts Table column widths in
# points, where 1 point = 12700 EMU.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

$emuPerPoint = 12700.0

$tbl.Columns.Item(1).Width = 1629100 / $emuPerPoint
$tbl.Columns.Item(2).Width = 5112912 / $emuPerPoint
$tbl.Columns.Item(3).Width = 1702079 / $emuPerPoint
